$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("M2").Value = 0.01569233333333333
$ws.Range("N2").Value = 0.047077
$ws.Range("O2").Value = 0.03693539111407157
$ws.Range("P2").Value = 0.03693539111407157
$ws.Range("Q2").Value = 0.005857576648111111
$ws.Range("R2").Value = 0.052718189833
$ws.Range("S2").Value = 0.01312779451718799
$ws.Range("T2").Value = 0.01312779451718799
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("N3").Value = 0.9690430000000001
$ws.Range("O3").Value = 0.7602859615386125
$ws.Range("P3").Value = 0.7602859615386125
$ws.Range("Q3").Value = 0.1205736059607778
$ws.Range("R3").Value = 1.085162453647
$ws.Range("S3").Value = 0.2702253198444973
$ws.Range("T3").Value = 0.2702253198444973
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3732763333333333
$ws.Range("H4").Value = 1.119829
$ws.Range("I4").Value = 0.3554258969843855
$ws.Range("J4").Value = 0.3554258969843855
$ws.Range("M4").Value = 0.08615233333333333
$ws.Range("N4").Value = 0.258457
$ws.Range("O4").Value = 0.202778647347316
$ws.Range("P4").Value = 0.202778647347316
$ws.Range("Q4").Value = 0.03215862709477778
$ws.Range("R4").Value = 0.289427643853
$ws.Range("S4").Value = 0.07207278262270017
$ws.Range("T4").Value = 0.07207278262270017
$ws.Range("G5").Value = 0.668317
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("M5").Value = 0.01569233333333333
$ws.Range("N5").Value = 0.047077
$ws.Range("O5").Value = 0.03693539111407157
$ws.Range("P5").Value = 0.03693539111407157
$ws.Range("Q5").Value = 0.01048745313633333
$ws.Range("R5").Value = 0.09438707822700002
$ws.Range("S5").Value = 0.02350411066781676
$ws.Range("T5").Value = 0.02350411066781676
$ws.Range("G6").Value = 0.668317
$ws.Range("I6").Value = 0.6363574327729865
$ws.Range("J6").Value = 0.6363574327729865
$ws.Range("N6").Value = 0.9690430000000001
$ws.Range("O6").Value = 0.7602859615386125
$ws.Range("P6").Value = 0.7602859615386125
$ws.Range("Q6").Value = 0.2158759702103334
$ws.Range("S6").Value = 0.483813622658053
$ws.Range("T6").Value = 0.483813622658053
$ws.Range("G7").Value = 0.668317
$ws.Range("I7").Value = 0.6363574327729865
$ws.Range("J7").Value = 0.6363574327729865
$ws.Range("M7").Value = 0.08615233333333333
$ws.Range("N7").Value = 0.258457
$ws.Range("O7").Value = 0.202778647347316
$ws.Range("P7").Value = 0.202778647347316
$ws.Range("Q7").Value = 0.05757706895633333
$ws.Range("R7").Value = 0.518193620607
$ws.Range("S7").Value = 0.1290396994471168
$ws.Range("T7").Value = 0.1290396994471168
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008629333333333334
$ws.Range("H8").Value = 0.025888
$ws.Range("I8").Value = 0.008216670242627913
$ws.Range("J8").Value = 0.008216670242627911
$ws.Range("M8").Value = 0.01569233333333333
$ws.Range("N8").Value = 0.047077
$ws.Range("O8").Value = 0.03693539111407157
$ws.Range("P8").Value = 0.03693539111407157
$ws.Range("Q8").Value = 0.0001354143751111111
$ws.Range("R8").Value = 0.001218729376
$ws.Range("S8").Value = 0.0003034859290668153
$ws.Range("T8").Value = 0.0003034859290668152
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008629333333333334
$ws.Range("H9").Value = 0.025888
$ws.Range("I9").Value = 0.008216670242627913
$ws.Range("J9").Value = 0.008216670242627911
$ws.Range("N9").Value = 0.9690430000000001
$ws.Range("O9").Value = 0.7602859615386125
$ws.Range("P9").Value = 0.7602859615386125
$ws.Range("Q9").Value = 0.002787398353777778
$ws.Range("R9").Value = 0.025086585184
$ws.Range("S9").Value = 0.006247019036062068
$ws.Range("T9").Value = 0.006247019036062066
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008629333333333334
$ws.Range("H10").Value = 0.025888
$ws.Range("I10").Value = 0.008216670242627913
$ws.Range("J10").Value = 0.008216670242627911
$ws.Range("M10").Value = 0.08615233333333333
$ws.Range("N10").Value = 0.258457
$ws.Range("O10").Value = 0.202778647347316
$ws.Range("P10").Value = 0.202778647347316
$ws.Range("Q10").Value = 0.0007434372017777779
$ws.Range("R10").Value = 0.006690934816
$ws.Range("S10").Value = 0.001666165277499031
$ws.Range("T10").Value = 0.001666165277499031
